$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force text storage for numeric-looking Price values so Excel
# does not auto-convert them to numbers (source data is text,
# matching the original inline-string cells).
$textForceCells = @('D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D16', 'D17', 'D20', 'D22', 'D23', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D37', 'D38', 'D40', 'D41', 'D42', 'D43', 'D44', 'D47', 'D49', 'D50', 'D51')
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '43.193.90'
$ws.Range('E2').Value = '  -4.88%  '
$ws.Range('D3').Value = '2.236.85'
$ws.Range('E3').Value = '  -5.84%  '
$ws.Range('D5').Value = '319.69'
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('D6').Value = '100.33'
$ws.Range('E6').Value = '  -7.84%  '
$ws.Range('D7').Value = '0.586'
$ws.Range('E7').Value = '  -8.45%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '0.563'
$ws.Range('E9').Value = '  -8.48%  '
$ws.Range('D10').Value = '36.88'
$ws.Range('E10').Value = '  -9.96%  '
$ws.Range('D11').Value = '54.48'
$ws.Range('E11').Value = '  -2.68%  '
$ws.Range('D12').Value = '0.0829'
$ws.Range('E12').Value = '  -9.71%  '
$ws.Range('D13').Value = '7.68'
$ws.Range('E13').Value = '  -10.15%  '
$ws.Range('D14').Value = '0.109'
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('D15').Value = '2.581.28'
$ws.Range('E15').Value = '  -5.69%  '
$ws.Range('D16').Value = '0.863'
$ws.Range('E16').Value = '  -12.44%  '
$ws.Range('D17').Value = '14.41'
$ws.Range('E17').Value = '  -6.83%  '
$ws.Range('D18').Value = '2.241.23'
$ws.Range('E18').Value = '  -5.44%  '
$ws.Range('D19').Value = '43.160.10'
$ws.Range('E19').Value = '  -4.89%  '
$ws.Range('D20').Value = '14.52'
$ws.Range('E20').Value = '  -9.08%  '
$ws.Range('D21').Value = '0.0₃0967'
$ws.Range('E21').Value = '  -9.14%  '
$ws.Range('D22').Value = '6.52'
$ws.Range('E22').Value = '  -11.17%  '
$ws.Range('D23').Value = '65.38'
$ws.Range('E23').Value = '  -10.84%  '
$ws.Range('E24').Value = '  -11.65%  '
$ws.Range('D25').Value = '237.04'
$ws.Range('E25').Value = '  -9.24%  '
$ws.Range('D26').Value = '2.15'
$ws.Range('E26').Value = '  -8.72%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('E28').Value = '  +2.25%  '
$ws.Range('D29').Value = '10.03'
$ws.Range('E29').Value = '  -10.18%  '
$ws.Range('D30').Value = '2.23'
$ws.Range('E30').Value = '  -2.68%  '
$ws.Range('D31').Value = '6.38'
$ws.Range('E31').Value = '  -16.57%  '
$ws.Range('D32').Value = '35.43'
$ws.Range('E32').Value = '  -4.98%  '
$ws.Range('D33').Value = '20.48'
$ws.Range('E33').Value = '  -8.45%  '
$ws.Range('D34').Value = '0.0872'
$ws.Range('E34').Value = '  -9.95%  '
$ws.Range('D35').Value = '153.10'
$ws.Range('E35').Value = '  -8.08%  '
$ws.Range('E36').Value = '  -4.68%  '
$ws.Range('D37').Value = '3.17'
$ws.Range('E37').Value = '  +8.17%  '
$ws.Range('D38').Value = '1.96'
$ws.Range('E38').Value = '  +3.63%  '
$ws.Range('E39').Value = '  -7.99%  '
$ws.Range('D40').Value = '4.43'
$ws.Range('E40').Value = '  -6.24%  '
$ws.Range('D41').Value = '0.103'
$ws.Range('E41').Value = '  -11.92%  '
$ws.Range('D42').Value = '3.67'
$ws.Range('E42').Value = '  -9.51%  '
$ws.Range('D43').Value = '0.0323'
$ws.Range('E43').Value = '  -8.99%  '
$ws.Range('D44').Value = '12.85'
$ws.Range('E44').Value = '  -1.94%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = '1.791.10'
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('D47').Value = '86.47'
$ws.Range('E47').Value = '  -12.19%  '
$ws.Range('E48').Value = '  -10.12%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').Value = '5.33'
$ws.Range('E49').Value = '  -10.68%  '
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').Value = '76.06'
$ws.Range('E50').Value = '  -9.42%  '
$ws.Range('D51').Value = '59.02'
$ws.Range('E51').Value = '  -16.21%  '
